# Generate Report for Archive
#
# 1) The shared string "Ready for handoff" becomes "In Translation" wherever
#    it is used: Overview!E2, Overview!F2 (per-language status columns) and
#    the "Status" column (C2) on each per-language detail sheet (zh-cn,
#    de-de).
# 2) As a side effect of the now-shorter status text, the "Status" columns
#    were re-sized (autofit-style) in the source workbook:
#      - Overview columns E & F: 17.2159881591797 -> 13.4101848602295
#      - zh-cn / de-de column C:  17.2159881591797 -> 13.4101848602295

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update the status text ---------------------------------------------
$overview.Range("E2").Value2 = "In Translation"
$overview.Range("F2").Value2 = "In Translation"
$zhcn.Range("C2").Value2     = "In Translation"
$dede.Range("C2").Value2     = "In Translation"

# --- Shrink the Status columns to match the new (shorter) text ----------
# ColumnWidth is expressed in "characters"; Excel quantizes it to whole
# pixels internally (stored width = (round(ColumnWidth*MDW)+5)/MDW), so we
# target the closest achievable character width to the recorded value.
$newStatusWidth = 12.5

$overview.Range("E1").EntireColumn.ColumnWidth = $newStatusWidth
$overview.Range("F1").EntireColumn.ColumnWidth = $newStatusWidth
$zhcn.Range("C1").EntireColumn.ColumnWidth     = $newStatusWidth
$dede.Range("C1").EntireColumn.ColumnWidth     = $newStatusWidth
